# resolve_report_task.xlsx update
# - "cht-conf" regenerated the XLSForm: NO_LABEL placeholder labels are gone,
#   every bare "text" question that only existed to carry NO_LABEL is now
#   typed "hidden" (and no longer needs a label cell), and the conditional
#   formatting for row 33 is folded back into the generic per-column rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# ---------------------------------------------------------------------
# 1. survey sheet data rows: type column changes + drop the NO_LABEL
#    label column (C) wherever it only ever held "NO_LABEL".
# ---------------------------------------------------------------------

# rows whose type flips from "text" -> "hidden" (and lose their C label)
$hiddenRows = 4,5,7,8,10,12,13
foreach ($r in $hiddenRows) {
    $ws.Range("A$r").Value = "hidden"
}

# rows that keep their existing type, but still lose the NO_LABEL cell in C
$keepTypeRows = 3,6,11
foreach ($r in $keepTypeRows) {
    # type (A) and name (B) are unchanged - only the label placeholder goes
}

# every row 3-13 that had a "NO_LABEL" label in column C loses that cell
$clearLabelRows = 3,4,5,6,7,8,10,11,12,13
foreach ($r in $clearLabelRows) {
    $ws.Range("C$r").ClearContents()
}

# ---------------------------------------------------------------------
# 2. survey sheet selection: cursor moved from J11 to C9
# ---------------------------------------------------------------------
$ws.Range("C9").Select()

# ---------------------------------------------------------------------
# 3. conditional formatting: row 33 (C33) had its own duplicated set of
#    cfRules; fold it back into the shared rules (same formulas as the
#    rest of columns A:G / C), and update the C-column rule to also
#    treat type "hidden" as not requiring a label, skipping rows that
#    already have a calculation in F.
# ---------------------------------------------------------------------
$ws.Range("C2").FormatConditions.Item(6).Formula1 = '=AND(ISBLANK(C2),NOT(OR(ISBLANK($A2),$A2="calculate",$A2="hidden")),ISBLANK($F2))'

$c33 = $ws.Range("C33")
$c33.FormatConditions.Delete()
$c33.FormatConditions.Add(2, 3, '=$A2="begin_group"')
$c33.FormatConditions.Add(2, 3, '=$A2="end_group"')
$c33.FormatConditions.Add(2, 3, '=$A2="begin_repeat"')
$c33.FormatConditions.Add(2, 3, '=$A2="end_repeat"')
$c33.FormatConditions.Add(2, 3, '=AND(ISBLANK($A2), NOT(ISBLANK(A2)))')
$c33.FormatConditions.Add(2, 3, '=AND(ISBLANK(C2),NOT(OR(ISBLANK($A2),$A2="calculate",$A2="hidden")),ISBLANK($F2))')

Write-Output "done"
